# UC003 - Cancelar Solicitação de Diária
# v1.1 -> v1.1.1
#
# The content of TC3's 2nd step and TC4's single step get swapped between
# the two test cases: TC3 becomes a single-step test case (using the text
# that used to live under TC4), while TC4 becomes a two-step test case
# (using the two steps that used to live under TC3). The row layout for
# TC4's header block (Test Case ID / Description / Precondition / column
# header) shifts up by one row to accommodate TC3 shrinking by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$msg205 = "SYSTEM Identifica que a solicitação de diária está em situação diferente de 'SOLICITADA PARA EMPENHO' ou 'SOLICITADA PARA PRESTAÇÃO DE CONTAS'.  Impede o cancelamento e exibe mensagem de erro (MSG205 - Solcitação de diária não pode ser cancelada) para o usuário."
$msg102 = "SYSTEM Exibe a mensagem (MSG102 - Confirmar cancelamento)"
$msg217 = "SYSTEM Identifica que o usuário não informou uma justificativa para o cancelamento. Não efetiva o cancelamento e exibe mensagem de erro (MSG217 - Necessário informar uma justificativa para o cancelamento de solicitações`t) para o usuário."
$chefeNaoInforma = "Chefe Não informa o motivo do cancelamento."
$chefeClica = "Chefe Clica em confirmar."
$chefeInforma = "Chefe Informa o motivo do cancelamento."

# --- Step 1: move TC4's header block (rows 31-34) up to rows 30-33 ---
# First unmerge the cells in this area so that writing a 6-column array
# of values lands in every column instead of just the merge anchor cell.
$ws.Range("B32:D32").UnMerge()
$ws.Range("B33:F33").UnMerge()

# Copy the full 6-column (A:F) content of each row up by one row.
$ws.Range("A30:F30").Value2 = $ws.Range("A31:F31").Value2
$ws.Range("A31:F31").Value2 = $ws.Range("A32:F32").Value2
$ws.Range("A32:F32").Value2 = $ws.Range("A33:F33").Value2
$ws.Range("A33:F33").Value2 = $ws.Range("A34:F34").Value2

# Re-create the merged cells for the Description/Precondition rows at
# their new locations (moved from rows 32/33 to rows 31/32).
$ws.Range("B31:D31").Merge()
$ws.Range("B32:F32").Merge()

# --- Step 2: TC4 now gets TWO steps (rows 34 and 35), reusing the text
# that used to belong to TC3's two steps. Row 34 previously held the
# column-header labels (now moved to row 33), so C34/E34/F34 must be
# cleared out instead of keeping those leftover header labels.
$ws.Range("A34").Value2 = 1
$ws.Range("B34").Value2 = $chefeNaoInforma
$ws.Range("C34").Value2 = ""
$ws.Range("D34").Value2 = $msg102
$ws.Range("E34").Value2 = ""
$ws.Range("F34").Value2 = ""

$ws.Range("A35").Value2 = 2
$ws.Range("B35").Value2 = $chefeClica
$ws.Range("C35").Value2 = ""
$ws.Range("D35").Value2 = $msg217
$ws.Range("E35").Value2 = ""
$ws.Range("F35").Value2 = ""

# --- Step 3: TC3 now gets a SINGLE step (row 27), reusing the text that
# used to belong to TC4's single step ---
$ws.Range("A27").Value2 = 1
$ws.Range("B27").Value2 = $chefeInforma
$ws.Range("D27").Value2 = $msg205

# --- Step 4: row 28 (TC3's former 2nd step) is no longer used; clear it
# completely (values and formatting) so it becomes a blank row again ---
$ws.Range("A28:F28").Clear()
